# Forgot to update PoS baselines
# Update column S ("Baseline") values on the Accuracy sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Accuracy")

$updates = @{
    2  = 0.2118068965517241
    3  = 0.1805458229957766
    7  = 0.2757242757242757
    10 = 0.216893039049236
    11 = 0.2738805263656158
    13 = 0.2939434534301629
    14 = 0.2832591683289857
    15 = 0.2672127950068266
    16 = 0.2971181376820615
    19 = 0.2493173885859726
    24 = 0.2075536661425302
    25 = 0.2711942992381789
    26 = 0.2733596891596734
    27 = 0.2498671219559148
    29 = 0.2504936941240743
}

foreach ($row in $updates.Keys) {
    $ws.Range("S$row").Value = $updates[$row]
}
